$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.483.62'
$ws.Range('E2').Value = '  -0.49%  '
# Row 3
$ws.Range('D3').Value = '3.502.63'
$ws.Range('E3').Value = '  -0.44%  '
# Row 4
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.26'
$ws.Range('E5').Value = '  -2.19%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '194.46'
$ws.Range('E6').Value = '  +0.84%  '
# Row 7
$ws.Range('E7').Value = '  -1.88%  '
# Row 8
$ws.Range('E8').Value = '  -0.01%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.199'
$ws.Range('E9').Value = '  -6.09%  '
# Row 10
$ws.Range('E10').Value = '  -3.30%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.96'
$ws.Range('E11').Value = '  -1.16%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000297'
$ws.Range('E12').Value = '  -3.98%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.41'
$ws.Range('E13').Value = '  -2.06%  '
# Row 14
$ws.Range('D14').Value = '4.060.52'
$ws.Range('E14').Value = '  -0.50%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '597.08'
$ws.Range('E15').Value = '  -3.90%  '
# Row 16
$ws.Range('D16').Value = '69.679.69'
$ws.Range('E16').Value = '  -0.23%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.91'
$ws.Range('E17').Value = '  -0.74%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.61'
$ws.Range('E18').Value = '  -0.27%  '
# Row 19
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.123'
$ws.Range('E19').Value = '  +2.14%  '
# Row 20
$ws.Range('B20').Value = 'WrappedEther'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D20').Value = '3.489.89'
$ws.Range('E20').Value = '  -0.77%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.980'
$ws.Range('E21').Value = '  -1.12%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.77'
$ws.Range('E22').Value = '  +3.40%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +4.16%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '101.97'
$ws.Range('E24').Value = '  -7.07%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.60'
$ws.Range('E25').Value = '  -2.53%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.09'
$ws.Range('E26').Value = '  -0.64%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.78'
$ws.Range('E27').Value = '  -2.20%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.45'
$ws.Range('E28').Value = '  -2.67%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.92'
$ws.Range('E29').Value = '  -3.81%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.24'
$ws.Range('E30').Value = '  +7.23%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.95'
$ws.Range('E31').Value = '  -0.60%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.26'
$ws.Range('E32').Value = '  -2.20%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').Value = '  -2.63%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.08'
$ws.Range('E34').Value = '  -0.75%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.20'
$ws.Range('E35').Value = '  +2.60%  '
# Row 36
$ws.Range('D36').Value = '3.745.16'
$ws.Range('E36').Value = '  +2.13%  '
# Row 37
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.11%  '
# Row 38
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0806'
$ws.Range('E38').Value = '  +3.20%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.65'
$ws.Range('E39').Value = '  -0.43%  '
# Row 40
$ws.Range('E40').Value = '  -1.95%  '
# Row 41
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '494.97'
$ws.Range('E41').Value = '  -4.17%  '
# Row 42
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '35.93'
$ws.Range('E42').Value = '  -2.20%  '
# Row 43
$ws.Range('E43').Value = '  -4.22%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0448'
$ws.Range('E44').Value = '  -5.09%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('E45').Value = '  -4.55%  '
# Row 46
$ws.Range('E46').Value = '  -3.29%  '
# Row 47
$ws.Range('E47').Value = '  -3.17%  '
# Row 48
$ws.Range('E48').Value = '  +0.19%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.38'
$ws.Range('E49').Value = '  -4.70%  '
# Row 50
$ws.Range('B50').Value = 'OceanProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').Value = '  +0.88%  '
# Row 51
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000242'
$ws.Range('E51').Value = '  +0.18%  '
